$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Minor text/typo corrections in existing address cells ---

# Brunei Malay: strip stray leading spaces
$ws.Range("C18").Value = "Limbang, Sarawak, Malaysia"
$ws.Range("F18").Value = "Miri, Sarawak, Malaysia"

# Greek Thrace Xoraxane Romane: add postal code
$ws.Range("C53").Value = "Drosero, Xanthi, Thrace, 671 00 Greece."

# Kumzari: drop "and Iran"
$ws.Range("B79").Value = "Musandam Peninsula, Oman"

# Mapudungun: strip stray leading space
$ws.Range("C91").Value = "La Pampa, Argentina"

# Setswana (South African): strip stray trailing space
$ws.Range("B113").Value = "Taung, South Africa"

# Southeastern Pashayi: updated place name (keep trailing non-breaking space)
$nbsp = [char]0x00A0
$ws.Range("B122").Value = "Kuz Kunar, Nangarhar, Afghanistan" + $nbsp

# --- Fix Pitjantjatjara publication link: new DOI, now a real hyperlink ---
$ws.Hyperlinks.Add($ws.Range("G107"), "https://doi.org/10.1017/S0025100314000073", [Type]::Missing, [Type]::Missing, "https://doi.org/10.1017/S0025100314000073") | Out-Null

# --- Append two new languages ---

# Row 158: Niuean
$ws.Range("A158").Value = "Niuean"
$ws.Range("B158").Value = "Niue"
$ws.Range("G158").Value = "https://doi.org/10.1017/S0025100317000500"
$ws.Range("H158").Value = "Illustration"

# Row 159: Hawaiian
$ws.Range("A159").Value = "Hawaiian"
$ws.Range("B159").Value = "Hawaii, USA"
$ws.Range("G159").Value = "https://doi.org/10.1017/S0025100316000438"
$ws.Range("H159").Value = "Illustration"
$ws.Range("I159").Value = "https://www.internationalphoneticassociation.org/sites/default/files/JIPArecordings/Hawaiian.zip"
